$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets sheet1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 278
$ws1.Range("F4").Value = 34
$ws1.Range("F5").Value = 3368
$ws1.Range("F6").Value = 2138
$ws1.Range("F7").Value = 408
$ws1.Range("F8").Value = 157
$ws1.Range("F9").Value = 41
$ws1.Range("F10").Value = 26
$ws1.Range("F11").Value = 1245
$ws1.Range("F12").Value = 225
$ws1.Range("F13").Value = 1428
$ws1.Range("F14").Value = 108

# Sheet "全部类型" (Worksheets sheet4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 278
$ws4.Range("F4").Value = 34
$ws4.Range("F5").Value = 3368
$ws4.Range("F6").Value = 2138
$ws4.Range("F7").Value = 408
$ws4.Range("F9").Value = 157
$ws4.Range("F10").Value = 41
$ws4.Range("F11").Value = 26
$ws4.Range("F14").Value = 1245
$ws4.Range("F15").Value = 225
$ws4.Range("F16").Value = 1428
$ws4.Range("F17").Value = 108
